$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 33314.332
$ws.Range("J93").Value = 33314.332
$ws.Range("L93").Value = 33314.332
$ws.Range("N93").Value = -38306.332
$ws.Range("H95").Value = 33299.332
$ws.Range("J95").Value = 33299.332
$ws.Range("L95").Value = 33299.332
$ws.Range("N95").Value = -38791.332
$ws.Range("H113").Value = 2130.4375
$ws.Range("I113").Value = 2113.375
$ws.Range("K113").Value = 2113.375
$ws.Range("M113").Value = 1140.625
$ws.Range("H120").Value = 48251.25
$ws.Range("J120").Value = 48251.25
$ws.Range("L120").Value = 48251.25
$ws.Range("N120").Value = -57927.25
$ws.Range("H124").Value = 53982
$ws.Range("J124").Value = 53982
$ws.Range("L124").Value = 53982
$ws.Range("N124").Value = -63802
$ws.Range("H126").Value = 47772
$ws.Range("J126").Value = 47772
$ws.Range("L126").Value = 47772
$ws.Range("N126").Value = -57652
$ws.Range("H128").Value = 54984
$ws.Range("J128").Value = 54984
$ws.Range("L128").Value = 54984
$ws.Range("N128").Value = -64944
$ws.Range("H130").Value = 54992
$ws.Range("J130").Value = 54992
$ws.Range("L130").Value = 54992
$ws.Range("N130").Value = -65032

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2127.2273
$ws.Range("I2").Value = 2143.6875
$ws.Range("K2").Value = 2143.6875
$ws.Range("M2").Value = -2030.6875
$ws.Range("H61").Value = 2079.762
$ws.Range("I61").Value = 1934.641
$ws.Range("K61").Value = 1934.641
$ws.Range("M61").Value = -1722.641
$ws.Range("H94").Value = 33315
$ws.Range("J94").Value = 33315
$ws.Range("L94").Value = 33315
$ws.Range("N94").Value = -35117
$ws.Range("H101").Value = 48187.6
$ws.Range("J101").Value = 48187.6
$ws.Range("L101").Value = 48187.6
$ws.Range("N101").Value = -54677.6
$ws.Range("H104").Value = 31075
$ws.Range("J104").Value = 31075
$ws.Range("L104").Value = 31075
$ws.Range("N104").Value = -38063
$ws.Range("H105").Value = 49370
$ws.Range("J105").Value = 49370
$ws.Range("L105").Value = 49370
$ws.Range("N105").Value = -56358
$ws.Range("H106").Value = 47935.2
$ws.Range("J106").Value = 47935.2
$ws.Range("L106").Value = 47935.2
$ws.Range("N106").Value = -50459.2
$ws.Range("H116").Value = 2127.2273
$ws.Range("I116").Value = 2143.6875
$ws.Range("K116").Value = 2143.6875
$ws.Range("M116").Value = 150.3125
$ws.Range("H119").Value = 48988
$ws.Range("J119").Value = 48988
$ws.Range("L119").Value = 48988
$ws.Range("N119").Value = -58664
$ws.Range("H123").Value = 45864.4
$ws.Range("J123").Value = 47330.5
$ws.Range("L123").Value = 47330.5
$ws.Range("N123").Value = -57130.5
$ws.Range("H136").Value = 2079.762
$ws.Range("I136").Value = 1934.641
$ws.Range("K136").Value = 5803.923000000001
$ws.Range("M136").Value = -3253.923000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2127.2273
$ws.Range("I3").Value = 2143.6875
$ws.Range("K3").Value = 2143.6875
$ws.Range("M3").Value = -2029.6875
$ws.Range("H6").Value = 19412.445
$ws.Range("J6").Value = 19412.445
$ws.Range("L6").Value = 19412.445
$ws.Range("N6").Value = -19638.445
$ws.Range("H95").Value = 44374
$ws.Range("J95").Value = 44374
$ws.Range("L95").Value = 44374
$ws.Range("N95").Value = -49866
$ws.Range("H103").Value = 250997.33
$ws.Range("J103").Value = 250997.33
$ws.Range("L103").Value = 250997.33
$ws.Range("N103").Value = -253341.33
$ws.Range("H122").Value = 40517.332
$ws.Range("J122").Value = 40517.332
$ws.Range("L122").Value = 40517.332
$ws.Range("N122").Value = -50317.332
$ws.Range("H126").Value = 50776
$ws.Range("J126").Value = 50776
$ws.Range("L126").Value = 50776
$ws.Range("N126").Value = -60656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1860.9
$ws.Range("I16").Value = 1621.7778
$ws.Range("J16").Value = 4013
$ws.Range("K16").Value = 1621.7778
$ws.Range("L16").Value = 4013
$ws.Range("M16").Value = -1334.7778
$ws.Range("N16").Value = -4587
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H28").Value = 34826.285
$ws.Range("J28").Value = 34826.285
$ws.Range("L28").Value = 34826.285
$ws.Range("N28").Value = -35316.285
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H31").Value = 4532.7163
$ws.Range("I31").Value = 1831.6857
$ws.Range("K31").Value = 1831.6857
$ws.Range("M31").Value = -1536.6857
$ws.Range("H34").Value = 4532.7163
$ws.Range("I34").Value = 1831.6857
$ws.Range("K34").Value = 1831.6857
$ws.Range("M34").Value = -1629.6857
$ws.Range("H43").Value = 49649
$ws.Range("J43").Value = 49649
$ws.Range("L43").Value = 49649
$ws.Range("N43").Value = -50017
$ws.Range("H58").Value = 2106.8867
$ws.Range("I58").Value = 1825.8605
$ws.Range("J58").Value = 3315.3
$ws.Range("K58").Value = 1825.8605
$ws.Range("L58").Value = 3315.3
$ws.Range("M58").Value = -1622.8605
$ws.Range("N58").Value = -3721.3
$ws.Range("H92").Value = 31494.572
$ws.Range("J92").Value = 31494.572
$ws.Range("L92").Value = 31494.572
$ws.Range("N92").Value = -36486.572
$ws.Range("H96").Value = 89082.664
$ws.Range("J96").Value = 89082.664
$ws.Range("L96").Value = 89082.664
$ws.Range("N96").Value = -94574.664
$ws.Range("H99").Value = 2608
$ws.Range("I99").Value = 2662
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2662
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -1164
$ws.Range("N99").Value = -5496
$ws.Range("H101").Value = 49649
$ws.Range("J101").Value = 49649
$ws.Range("L101").Value = 49649
$ws.Range("N101").Value = -56139
$ws.Range("H106").Value = 42994.4
$ws.Range("J106").Value = 42994.4
$ws.Range("L106").Value = 42994.4
$ws.Range("N106").Value = -45518.4
$ws.Range("H112").Value = 40997.332
$ws.Range("J112").Value = 40997.332
$ws.Range("L112").Value = 40997.332
$ws.Range("N112").Value = -43951.332
$ws.Range("H113").Value = 1860.9
$ws.Range("I113").Value = 1621.7778
$ws.Range("J113").Value = 4013
$ws.Range("K113").Value = 1621.7778
$ws.Range("L113").Value = 4013
$ws.Range("M113").Value = 548.2221999999999
$ws.Range("N113").Value = -8353
$ws.Range("H118").Value = 44742
$ws.Range("J118").Value = 44742
$ws.Range("L118").Value = 44742
$ws.Range("N118").Value = -48056
$ws.Range("H122").Value = 111247600
$ws.Range("I122").Value = 166867570
$ws.Range("J122").Value = 7671.3335
$ws.Range("K122").Value = 500602710
$ws.Range("L122").Value = 23014.0005
$ws.Range("M122").Value = -500600260
$ws.Range("N122").Value = -27914.0005
$ws.Range("H124").Value = 45326
$ws.Range("J124").Value = 45326
$ws.Range("L124").Value = 45326
$ws.Range("N124").Value = -50236
$ws.Range("H126").Value = 2608
$ws.Range("I126").Value = 2662
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 7986
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -5516
$ws.Range("N126").Value = -12440
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 38326
$ws.Range("J131").Value = 38326
$ws.Range("L131").Value = 38326
$ws.Range("N131").Value = -48406
$ws.Range("H132").Value = 80976
$ws.Range("I132").Value = 2694.75
$ws.Range("J132").Value = 143601
$ws.Range("K132").Value = 8084.25
$ws.Range("L132").Value = 430803
$ws.Range("M132").Value = -5554.25
$ws.Range("N132").Value = -435863
$ws.Range("H136").Value = 2106.8867
$ws.Range("I136").Value = 1825.8605
$ws.Range("J136").Value = 3315.3
$ws.Range("K136").Value = 5477.5815
$ws.Range("L136").Value = 9945.900000000001
$ws.Range("M136").Value = -2927.5815
$ws.Range("N136").Value = -15045.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 7337.067
$ws.Range("I113").Value = 8327.385
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 24982.155
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = -22812.155
$ws.Range("N113").Value = -7040
$ws.Range("H137").Value = 52641264
$ws.Range("I137").Value = 4288.4287
$ws.Range("J137").Value = 83346170
$ws.Range("K137").Value = 12865.2861
$ws.Range("L137").Value = 250038510
$ws.Range("M137").Value = -7765.286100000001
$ws.Range("N137").Value = -250048710
$ws.Range("H139").Value = 7525.8335
$ws.Range("I139").Value = 16526.25
$ws.Range("J139").Value = 3025.625
$ws.Range("K139").Value = 49578.75
$ws.Range("L139").Value = 9076.875
$ws.Range("M139").Value = -44438.75
$ws.Range("N139").Value = -19356.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 195942.23
$ws.Range("I80").Value = 459018.2
$ws.Range("J80").Value = 3019.8667
$ws.Range("K80").Value = 459018.2
$ws.Range("L80").Value = 3019.8667
$ws.Range("M80").Value = -458020.2
$ws.Range("N80").Value = -5015.8667
$ws.Range("H83").Value = 195942.23
$ws.Range("I83").Value = 459018.2
$ws.Range("J83").Value = 3019.8667
$ws.Range("K83").Value = 2295091
$ws.Range("L83").Value = 15099.3335
$ws.Range("M83").Value = -2290099
$ws.Range("N83").Value = -25083.3335
$ws.Range("H96").Value = 35650.332
$ws.Range("J96").Value = 35650.332
$ws.Range("L96").Value = 35650.332
$ws.Range("N96").Value = -41142.332
$ws.Range("H98").Value = 50023700
$ws.Range("J98").Value = 50023700
$ws.Range("L98").Value = 50023700
$ws.Range("N98").Value = -50029690
$ws.Range("H100").Value = 37340
$ws.Range("J100").Value = 37340
$ws.Range("L100").Value = 37340
$ws.Range("N100").Value = -39504
$ws.Range("H101").Value = 50657
$ws.Range("J101").Value = 50657
$ws.Range("L101").Value = 50657
$ws.Range("N101").Value = -57147
$ws.Range("H104").Value = 42890.332
$ws.Range("J104").Value = 42890.332
$ws.Range("L104").Value = 42890.332
$ws.Range("N104").Value = -49878.332
$ws.Range("H105").Value = 47671
$ws.Range("J105").Value = 47671
$ws.Range("L105").Value = 47671
$ws.Range("N105").Value = -54659
$ws.Range("H113").Value = 2620
$ws.Range("I113").Value = 2433.3333
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 2433.3333
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -263.3332999999998
$ws.Range("N113").Value = -7240
$ws.Range("H116").Value = 49538
$ws.Range("J116").Value = 49538
$ws.Range("L116").Value = 49538
$ws.Range("N116").Value = -58716
$ws.Range("H118").Value = 36204
$ws.Range("J118").Value = 36204
$ws.Range("L118").Value = 36204
$ws.Range("N118").Value = -39518
$ws.Range("H127").Value = 23660.25
$ws.Range("J127").Value = 23660.25
$ws.Range("L127").Value = 23660.25
$ws.Range("N127").Value = -33580.25
$ws.Range("H128").Value = 38514.668
$ws.Range("J128").Value = 38514.668
$ws.Range("L128").Value = 38514.668
$ws.Range("N128").Value = -48474.668
$ws.Range("H131").Value = 38986
$ws.Range("J131").Value = 38986
$ws.Range("L131").Value = 38986
$ws.Range("N131").Value = -49066

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 45737.625
$ws.Range("J94").Value = 45737.625
$ws.Range("L94").Value = 45737.625
$ws.Range("N94").Value = -47089.625
$ws.Range("H98").Value = 44245
$ws.Range("J98").Value = 44245
$ws.Range("L98").Value = 44245
$ws.Range("N98").Value = -50235
$ws.Range("H108").Value = 48626
$ws.Range("J108").Value = 48626
$ws.Range("L108").Value = 48626
$ws.Range("N108").Value = -56306
$ws.Range("H109").Value = 28325.334
$ws.Range("J109").Value = 28325.334
$ws.Range("L109").Value = 28325.334
$ws.Range("N109").Value = -31099.334
$ws.Range("H111").Value = 44383
$ws.Range("J111").Value = 44383
$ws.Range("L111").Value = 44383
$ws.Range("N111").Value = -52563
$ws.Range("H112").Value = 43387
$ws.Range("J112").Value = 43387
$ws.Range("L112").Value = 43387
$ws.Range("N112").Value = -46341
$ws.Range("H114").Value = 38136.5
$ws.Range("J114").Value = 38136.5
$ws.Range("L114").Value = 38136.5
$ws.Range("N114").Value = -46814.5
$ws.Range("H116").Value = 47286
$ws.Range("J116").Value = 47286
$ws.Range("L116").Value = 47286
$ws.Range("N116").Value = -56464
$ws.Range("H117").Value = 25996
$ws.Range("J117").Value = 25996
$ws.Range("L117").Value = 25996
$ws.Range("N117").Value = -35174
$ws.Range("H120").Value = 56648.668
$ws.Range("J120").Value = 56648.668
$ws.Range("L120").Value = 56648.668
$ws.Range("N120").Value = -66324.66800000001
$ws.Range("H121").Value = 31524.2
$ws.Range("J121").Value = 31524.2
$ws.Range("L121").Value = 31524.2
$ws.Range("N121").Value = -35018.2
$ws.Range("H123").Value = 29781.3
$ws.Range("J123").Value = 29781.3
$ws.Range("L123").Value = 29781.3
$ws.Range("N123").Value = -39581.3
$ws.Range("H131").Value = 32155
$ws.Range("J131").Value = 32155
$ws.Range("L131").Value = 32155
$ws.Range("N131").Value = -42235
$ws.Range("H132").Value = 5412.48
$ws.Range("I132").Value = 6354.364
$ws.Range("J132").Value = 4672.4287
$ws.Range("K132").Value = 19063.092
$ws.Range("L132").Value = 14017.2861
$ws.Range("M132").Value = -16533.092
$ws.Range("N132").Value = -19077.2861
$ws.Range("H136").Value = 2080.1292
$ws.Range("I136").Value = 1788.56
$ws.Range("J136").Value = 3295
$ws.Range("K136").Value = 5365.68
$ws.Range("L136").Value = 9885
$ws.Range("M136").Value = -2815.68
$ws.Range("N136").Value = -14985
$ws.Range("H140").Value = 38931.668
$ws.Range("J140").Value = 38931.668
$ws.Range("L140").Value = 38931.668
$ws.Range("N140").Value = -49291.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11514
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 11514
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 11514
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -13386
$ws.Range("H77").Value = 11514
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 11514
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 34542
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -43902
$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992
$ws.Range("H94").Value = 21000
$ws.Range("J94").Value = 21000
$ws.Range("L94").Value = 21000
$ws.Range("N94").Value = -22802
$ws.Range("H95").Value = 40781.332
$ws.Range("J95").Value = 40781.332
$ws.Range("L95").Value = 40781.332
$ws.Range("N95").Value = -46273.332
$ws.Range("H97").Value = 38054.668
$ws.Range("J97").Value = 38054.668
$ws.Range("L97").Value = 38054.668
$ws.Range("N97").Value = -40036.668
$ws.Range("H98").Value = 45303
$ws.Range("J98").Value = 45303
$ws.Range("L98").Value = 45303
$ws.Range("N98").Value = -51293
$ws.Range("H102").Value = 37456
$ws.Range("J102").Value = 37456
$ws.Range("L102").Value = 37456
$ws.Range("N102").Value = -43946
$ws.Range("H103").Value = 36963
$ws.Range("J103").Value = 36963
$ws.Range("L103").Value = 36963
$ws.Range("N103").Value = -39307
$ws.Range("H108").Value = 48626
$ws.Range("J108").Value = 48626
$ws.Range("L108").Value = 48626
$ws.Range("N108").Value = -56306
$ws.Range("H109").Value = 36871
$ws.Range("J109").Value = 36871
$ws.Range("L109").Value = 36871
$ws.Range("N109").Value = -39645
$ws.Range("H112").Value = 36387
$ws.Range("J112").Value = 36387
$ws.Range("L112").Value = 36387
$ws.Range("N112").Value = -39341
$ws.Range("H114").Value = 34380
$ws.Range("J114").Value = 34380
$ws.Range("L114").Value = 34380
$ws.Range("N114").Value = -43058
$ws.Range("H115").Value = 37366.332
$ws.Range("J115").Value = 37366.332
$ws.Range("L115").Value = 37366.332
$ws.Range("N115").Value = -40500.332
$ws.Range("H116").Value = 48690.668
$ws.Range("J116").Value = 48690.668
$ws.Range("L116").Value = 48690.668
$ws.Range("N116").Value = -57868.668
$ws.Range("H118").Value = 38171.668
$ws.Range("J118").Value = 38171.668
$ws.Range("L118").Value = 38171.668
$ws.Range("N118").Value = -41485.668
$ws.Range("H119").Value = 48998.4
$ws.Range("J119").Value = 48998.4
$ws.Range("L119").Value = 48998.4
$ws.Range("N119").Value = -58674.4
$ws.Range("H120").Value = 43312.4
$ws.Range("J120").Value = 43312.4
$ws.Range("L120").Value = 43312.4
$ws.Range("N120").Value = -52988.4
$ws.Range("H121").Value = 43394.668
$ws.Range("J121").Value = 43394.668
$ws.Range("L121").Value = 43394.668
$ws.Range("N121").Value = -46888.668
$ws.Range("H125").Value = 39715
$ws.Range("J125").Value = 39715
$ws.Range("L125").Value = 39715
$ws.Range("N125").Value = -49555
$ws.Range("H129").Value = 42406.4
$ws.Range("J129").Value = 42406.4
$ws.Range("L129").Value = 42406.4
$ws.Range("N129").Value = -52406.4
$ws.Range("H132").Value = 1426
$ws.Range("I132").Value = 951.9545000000001
$ws.Range("K132").Value = 2855.8635
$ws.Range("M132").Value = -325.8635000000004
